# Auto-generated edit script: updates Universalis market-price-derived
# columns (H:N) on the per-job leve-profit sheets (ALC, ARM, BSM, CRP,
# CUL, LTW, WVR) to reflect a refreshed market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 623.86664
$ws.Range("I6").Value = 33.11111
$ws.Range("J6").Value = 1510
$ws.Range("K6").Value = 99.33332999999999
$ws.Range("L6").Value = 4530
$ws.Range("M6").Value = 12.66667000000001
$ws.Range("N6").Value = -4754
$ws.Range("H15").Value = 1428.537
$ws.Range("I15").Value = 1428.537
$ws.Range("K15").Value = 4285.611
$ws.Range("M15").Value = -4116.611
$ws.Range("H17").Value = 1910.0294
$ws.Range("I17").Value = 1444.4445
$ws.Range("J17").Value = 2077.64
$ws.Range("K17").Value = 4333.333500000001
$ws.Range("L17").Value = 6232.92
$ws.Range("M17").Value = -4165.333500000001
$ws.Range("N17").Value = -6568.92
$ws.Range("H58").Value = 1664.7858
$ws.Range("I58").Value = 233.66667
$ws.Range("J58").Value = 2738.125
$ws.Range("K58").Value = 701.00001
$ws.Range("L58").Value = 8214.375
$ws.Range("M58").Value = -551.00001
$ws.Range("N58").Value = -8514.375
$ws.Range("H74").Value = 9928.571
$ws.Range("I74").Value = 3900
$ws.Range("K74").Value = 3900
$ws.Range("M74").Value = -2964
$ws.Range("H77").Value = 9928.571
$ws.Range("I77").Value = 3900
$ws.Range("K77").Value = 19500
$ws.Range("M77").Value = -14820
$ws.Range("H116").Value = 3642.5
$ws.Range("I116").Value = 3642.5
$ws.Range("K116").Value = 3642.5
$ws.Range("M116").Value = -200.5
$ws.Range("H131").Value = 1235.7
$ws.Range("I131").Value = 844.1667
$ws.Range("K131").Value = 2532.5001
$ws.Range("M131").Value = 2507.4999
$ws.Range("H132").Value = 19131.166
$ws.Range("I132").Value = 24708.223
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 74124.66900000001
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -71594.66900000001
$ws.Range("N132").Value = -12260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3738.9167
$ws.Range("I45").Value = 2552.8572
$ws.Range("K45").Value = 2552.8572
$ws.Range("M45").Value = -2175.8572
$ws.Range("H86").Value = 10642.5
$ws.Range("I86").Value = 6285
$ws.Range("J86").Value = 15000
$ws.Range("K86").Value = 6285
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = -5099
$ws.Range("N86").Value = -17372
$ws.Range("H89").Value = 10642.5
$ws.Range("I89").Value = 6285
$ws.Range("J89").Value = 15000
$ws.Range("K89").Value = 18855
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -12927
$ws.Range("N89").Value = -56856
$ws.Range("H97").Value = 1371.4
$ws.Range("I97").Value = 930
$ws.Range("K97").Value = 930
$ws.Range("M97").Value = -434
$ws.Range("H102").Value = 8214.286
$ws.Range("I102").Value = 6875
$ws.Range("K102").Value = 6875
$ws.Range("M102").Value = -5253

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9333.166999999999
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 9999.799999999999
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 9999.799999999999
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -12245.8
$ws.Range("H89").Value = 9333.166999999999
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 9999.799999999999
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 49999
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -61231
$ws.Range("H99").Value = 5869.6665
$ws.Range("I99").Value = 5554.75
$ws.Range("J99").Value = 6499.5
$ws.Range("K99").Value = 5554.75
$ws.Range("L99").Value = 6499.5
$ws.Range("M99").Value = -4056.75
$ws.Range("N99").Value = -9495.5
$ws.Range("H132").Value = 155000
$ws.Range("J132").Value = 155000
$ws.Range("L132").Value = 155000
$ws.Range("N132").Value = -165120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2827.9473
$ws.Range("I58").Value = 1880.3846
$ws.Range("K58").Value = 1880.3846
$ws.Range("M58").Value = -1677.3846
$ws.Range("H95").Value = 17569.455
$ws.Range("J95").Value = 17569.455
$ws.Range("L95").Value = 17569.455
$ws.Range("N95").Value = -23061.455
$ws.Range("H132").Value = 4984.5
$ws.Range("I132").Value = 4210.4443
$ws.Range("J132").Value = 6377.8
$ws.Range("K132").Value = 12631.3329
$ws.Range("L132").Value = 19133.4
$ws.Range("M132").Value = -10101.3329
$ws.Range("N132").Value = -24193.4
$ws.Range("H136").Value = 2827.9473
$ws.Range("I136").Value = 1880.3846
$ws.Range("K136").Value = 5641.1538
$ws.Range("M136").Value = -3091.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 335.66666
$ws.Range("J38").Value = 2
$ws.Range("L38").Value = 6
$ws.Range("N38").Value = -700
$ws.Range("H59").Value = 800
$ws.Range("I59").Value = 800
$ws.Range("K59").Value = 2400
$ws.Range("M59").Value = -1860
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11685
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -10908
$ws.Range("H131").Value = 2859.6
$ws.Range("J131").Value = 2928.4285
$ws.Range("L131").Value = 8785.2855
$ws.Range("N131").Value = -18865.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5535.478
$ws.Range("I7").Value = 4240
$ws.Range("J7").Value = 7964.5
$ws.Range("K7").Value = 4240
$ws.Range("L7").Value = 7964.5
$ws.Range("M7").Value = -4128
$ws.Range("N7").Value = -8188.5
$ws.Range("H16").Value = 999.25
$ws.Range("I16").Value = 999.25
$ws.Range("K16").Value = 999.25
$ws.Range("M16").Value = -829.25
$ws.Range("H22").Value = 2179.8
$ws.Range("J22").Value = 2333
$ws.Range("L22").Value = 2333
$ws.Range("N22").Value = -2923
$ws.Range("H27").Value = 2179.8
$ws.Range("J27").Value = 2333
$ws.Range("L27").Value = 2333
$ws.Range("N27").Value = -2547
$ws.Range("H46").Value = 4924.8945
$ws.Range("I46").Value = 3382.3333
$ws.Range("J46").Value = 5636.846
$ws.Range("K46").Value = 3382.3333
$ws.Range("L46").Value = 5636.846
$ws.Range("M46").Value = -3194.3333
$ws.Range("N46").Value = -6012.846
$ws.Range("H93").Value = 1662.6154
$ws.Range("I93").Value = 1851.5
$ws.Range("K93").Value = 1851.5
$ws.Range("M93").Value = -603.5
$ws.Range("H100").Value = 6175.9375
$ws.Range("I100").Value = 3057.077
$ws.Range("J100").Value = 8309.895
$ws.Range("K100").Value = 3057.077
$ws.Range("L100").Value = 8309.895
$ws.Range("M100").Value = -2516.077
$ws.Range("N100").Value = -9391.895
$ws.Range("H126").Value = 5535.478
$ws.Range("I126").Value = 4240
$ws.Range("J126").Value = 7964.5
$ws.Range("K126").Value = 12720
$ws.Range("L126").Value = 23893.5
$ws.Range("M126").Value = -10250
$ws.Range("N126").Value = -28833.5
$ws.Range("H132").Value = 6287.8
$ws.Range("I132").Value = 4975
$ws.Range("J132").Value = 7163
$ws.Range("K132").Value = 14925
$ws.Range("L132").Value = 21489
$ws.Range("M132").Value = -12395
$ws.Range("N132").Value = -26549
$ws.Range("H136").Value = 4084.4285
$ws.Range("I136").Value = 4084.4285
$ws.Range("K136").Value = 12253.2855
$ws.Range("M136").Value = -9703.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 55683.332
$ws.Range("J39").Value = 55683.332
$ws.Range("L39").Value = 55683.332
$ws.Range("N39").Value = -56509.332
$ws.Range("H132").Value = 2591.2354
$ws.Range("J132").Value = 3125
$ws.Range("L132").Value = 9375
$ws.Range("N132").Value = -14435
